$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 4.834600000000002
$ws.Range("D4").Value = -7.804000000000005

$ws.Range("D5").Value = -8.102999999999998

$ws.Range("B7").Value = 6.328199999999994

$ws.Range("D8").Value = -8.171599999999998

$ws.Range("B16").Value = 8.404800000000007
$ws.Range("D16").Value = -8.139600000000002
